# SAGM_PROCESOS.xlsx edit script
# - Adds two new rows to Sheet1 documenting Chart.js installation steps
# - Moves the active/selected tab from "Reports" (sheet3) back to "Sheet1"
# - Updates selection/active-cell bookkeeping to match

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item("Sheet1")
$sheet2 = $wb.Worksheets.Item("MIGRACIONES DE FRAMEWORK")
$sheet3 = $wb.Worksheets.Item("Reports")

# New content rows on Sheet1 (C13:D14)
$sheet1.Range("C13").Value = "INSTALACION Chart.js"
$sheet1.Range("D13").Value = 'se instala desde consola de NPM y despues vas a "C:\Users\manol\.nuget\packages\chart.js\3.7.1\content\Scripts"'
$sheet1.Range("D14").Value = "y copiamos la librería"

# Update selections to match the saved view state
$sheet2.Activate()
$sheet2.Range("C4").Select()
$excel.ActiveWindow.ScrollRow = 10      # scroll so row 10 becomes the top-left visible cell

$sheet3.Activate()
$sheet3.Range("D10").Select()

# Make Sheet1 the active/selected sheet (was "Reports") and restore its selection
$sheet1.Activate()
$sheet1.Range("C12").Select()
$sheet1.Select()
